# Rename the inline picture shapes that live in the document's headers
# and footers:
#   - the two "PearsonLogo.png" pictures (in the primary and first-page
#     footers) go from "image2.png" -> "image1.png"
#   - the "BTec_Logo-Orange" picture (in the first-page header) goes
#     from "image1.jpg" -> "image2.jpg"
#
# Word exposes the underlying OOXML <wp:docPr>/<pic:cNvPr> "name" as the
# single InlineShape.Name property, so renaming is just a property
# assignment. Re-fetching the InlineShape through its own (single
# character) Range before writing the Name keeps the object handle from
# going stale when the shape is not the very first thing in a
# multi-paragraph header/footer story.

function Set-InlineShapeName {
    param($shape, [string]$newName)

    # Re-seat the handle via its own Range so the write always targets a
    # freshly addressed block, then assign the new name.
    $fresh = $shape.Range.InlineShapes.Item(1)
    $fresh.Name = $newName
}

$d = $word.ActiveDocument
$sec = $d.Sections.First

# Walk every header/footer in the (only) section and rename whichever
# picture we find by matching its alt text / description, so the script
# does not depend on a particular Headers/Footers index ordering.
$targets = @()
for ($i = 1; $i -le 3; $i++) {
    $h = $sec.Headers.Item($i)
    if ($h.Exists) { $targets += $h }
}
for ($i = 1; $i -le 3; $i++) {
    $f = $sec.Footers.Item($i)
    if ($f.Exists) { $targets += $f }
}

foreach ($hf in $targets) {
    $shapes = $hf.Range.InlineShapes
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        $descr = $shp.AlternativeText

        if ($descr -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
            Set-InlineShapeName $shp "image1.png"
        }
        elseif ($descr -eq "BTec_Logo-Orange") {
            Set-InlineShapeName $shp "image2.jpg"
        }
    }
}
